$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Text"
$ws.Range("B2:B5").Value = "causalclaims"
$ws.Range("B6").Select()
